$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dataset")

# Insert a new column before column I (species moves from I to J, etc.)
$ws.Columns.Item(9).Insert()

# Set the header text for the newly inserted column I
$ws.Range("I1").Value = "has_related_ids"
